# The author replaced the placeholder text in B1 ("Metadaten aus Excel")
# with a long German paragraph describing Lower Saxony / Germany population
# statistics, turned word-wrap on for that cell, and grew row 1 to fit it
# (Excel's maximum row height, 409.5pt). A2/B2 keep their original content
# ("other_info" / "Weitere Informationen") - the shared-string table just
# gets re-ordered as a side effect of removing the old string and appending
# the new long one at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$longText = "Anstieg der Bevölkerung deutlich höher als im Vorjahr Die Bevölkerungszahl Niedersachsens lag auch am Ende des Jahres 2021 bei über 8 Mio. Einwohnerinnen und Einwohnern. Im Vergleich zum Vorjahr entsprach dies einem Plus von 23 610 Personen bzw. rund 0,3 %. Der Anstieg fiel damit größer aus als im Vorjahr (2019/2020: +9 813 Personen). In Deutschland lebten Ende 2021 insgesamt über 83,2 Mio. Menschen; die Einwohnerzahl Deutschlands stieg somit nach einem leichten Rückgang im Vorjahr wieder um 82 093 Personen bzw. um 0,1 %. Die Bevölkerungsentwicklung war 2021 in Niedersachsen damit im zehnten Jahr in Folge`nsteigend, ebenso wie in 6 weiteren Ländern. Währenddessen gab es das fünfte Jahr in`nFolge eine rückgängige Bevölkerungsentwicklung in 4 Ländern: Im Vergleich 2021 gegenüber 2016 betrug der Rückgang in Sachsen-Anhalt 3,0 %, in Thüringen 2,3 %, im Saarland 1,4 % und in Sachsen 1,0 %. In diesem Zeitraum stieg die Zahl der Einwohnerinnen`nund Einwohner in Niedersachsen um rund 81 300 Personen bzw. 1,0 %. In Deutschland`nwaren es etwas über 700 000 Personen mehr (+0,9 %).`nDurchschnittlich wohnten 2021 in Niedersachsen 168 Einwohnerinnen und Einwohner je`nQuadratkilometer (Deutschland: 233 Personen). Im Ländervergleich hatten nur 4 Länder`neine geringere Bevölkerungsdichte: Mecklenburg-Vorpommern (69), Brandenburg (86),`nSachsen-Anhalt (106) und Thüringen (130)."

# Replace B1's short placeholder with the long paragraph.
$ws.Range("B1").Value = $longText

# Wrap the text and grow row 1 so the whole paragraph is visible
# (matches Excel's own max row height cap of 409.5 points).
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 409.5

# Reflect the author's final on-screen selection/scroll state.
$ws.Range("B1").Select() | Out-Null
